$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("team")

# Insert two new rows above the "Franco Muriñigo" row (currently row 11),
# pushing it and "Thatiane Pereira" down to rows 13 and 14.
$ws.Range("A11:A12").EntireRow.Insert()

# The insert copies formatting down from the row above (row 10, which has
# hyperlink-style cells in columns E/F); strip that back to the sheet
# default before writing the new values.
$ws.Range("A11:F12").ClearFormats()

# New team member: Caetano Rocha
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Caetano Rocha"
$ws.Range("C11").Value = "Visiting Scholar"
$ws.Range("D11").Value = "caetano_rocha.jpg"
$ws.Range("E11").Value = "https://www.linkedin.com/in/caetano-rocha-9087b0319/"
$ws.Range("F11").Value = "Agronomy undergraduate from University of Santa Maria with experience in field crops and soil science. Supports field operations through crop monitoring, data collection, and hands-on agricultural practices across research sites."

# New team member: Diogo Verzegnazzi
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Diogo Verzegnazzi"
$ws.Range("C12").Value = "Visiting Scholar"
$ws.Range("D12").Value = "diogo_verzegnazzi.jpg"
$ws.Range("E12").Value = "https://www.linkedin.com/in/diogo-verzegnazzi-555418239/"
$ws.Range("F12").Value = "Agronomy undergraduate from University of Santa Maria with experience in field crops. Supports field operations through crop monitoring, data collection, and hands-on agricultural practices across research sites."

# Renumber the "id" column for the rows that shifted down.
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13

$ws.Range("E17").Select()
